$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Week 4 odds/team rows appended below the existing Week 1-3 data.
# Column A = Odds, Column B = Team, Column C = Week number.
$data = @(
    @(320, 'Georgia', 4),
    @(350, 'Ohio State', 4),
    @(500, 'Texas', 4),
    @(900, 'Alabama', 4),
    @(1100, 'Oregon', 4),
    @(1200, 'Ole Miss', 4),
    @(1600, 'Tennessee', 4),
    @(1800, 'Miami', 4),
    @(2000, 'Penn State', 4),
    @(3000, 'USC', 4),
    @(4000, 'LSU', 4),
    @(5000, 'Utah', 4),
    @(4000, 'Missouri', 4),
    @(4000, 'Notre Dame', 4),
    @(6000, 'Clemson', 4),
    @(6000, 'Kansas State', 4),
    @(8000, 'Oklahoma', 4),
    @(8000, 'Texas A&M', 4),
    @(12500, 'Michigan', 4),
    @(12500, 'Louisville', 4),
    @(15000, 'Nebraska', 4),
    @(15000, 'Oklahoma State', 4),
    @(50000, 'Iowa', 4),
    @(40000, 'South Carolina', 4),
    @(30000, 'Auburn', 4),
    @(50000, 'Wisconsin', 4),
    @(50000, 'Texas Tech', 4),
    @(30000, 'Boise State', 4),
    @(75000, 'Oregon State', 4),
    @(50000, 'Colorado', 4),
    @(20000, 'Iowa State', 4),
    @(30000, 'Syracuse', 4),
    @(30000, 'UCF', 4),
    @(50000, 'Kentucky', 4),
    @(50000, 'Washington', 4),
    @(50000, 'West Virginia', 4),
    @(40000, 'Arizona', 4),
    @(30000, 'Memphis', 4),
    @(50000, 'Florida', 4),
    @(100000, 'UNLV', 4),
    @(30000, 'NC State', 4),
    @(40000, 'Kansas', 4),
    @(100000, 'Washington State', 4),
    @(100000, 'Baylor', 4),
    @(100000, 'UCLA', 4),
    @(40000, 'Virginia Tech', 4),
    @(40000, 'TCU', 4),
    @(75000, 'SMU', 4),
    @(100000, 'Michigan State', 4),
    @(50000, 'Maryland', 4)
)

$startRow = 152
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Match the author's final view state: scrolled down near the bottom of the
# newly-added rows with cell C200 selected.
$ws.Activate()
$ws.Range("C200").Select()
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
